$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, reusing G1's formatting (bold font,
# border, centered/top alignment) so the new column matches the existing
# header row style exactly.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Per-row "Save" flag values (rows 2-48), taken from the era data update.
$saveValues = @(0,0,1,0,0,1,0,0,0,1,0,0,0,0,1,1,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,1,1,0,1,0,0,0,0,0,0,1,0,1,1,0,0)

$row = 2
foreach ($val in $saveValues) {
    $ws.Cells.Item($row, 8).Value = $val
    $row++
}
